# "turn negative links red"
#
# On slide 3, the curved connector representing the negative correlation
# (–0.49***) between "Epidemic size" and the other variable is recolored
# red, nudged slightly to the right, and its label text is recolored red
# to match.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# Shape 4 = "Connector: Curved 75" (the arrow pointing at the -0.49*** link)
$conn = $s.Shapes.Item(4)
$conn.Left = 317.6123
$conn.Line.ForeColor.RGB = 192

# Shape 5 = "TextBox 87" (contains the "-0.49***" label)
$lbl = $s.Shapes.Item(5)
$lbl.TextFrame.TextRange.Font.Color.RGB = 192
